$wb = $excel.ActiveWorkbook

# --- Status text change: "Ready for handoff" -> "In Translation" ---
# Every cell across the three sheets that shows this status needs updating
# so the shared string is fully replaced (no more references remain to the
# old text).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = "In Translation"
$wsOverview.Cells.Item(2, 6).Value = "In Translation"
$wsOverview.Cells.Item(3, 5).Value = "In Translation"
$wsOverview.Cells.Item(3, 6).Value = "In Translation"
$wsOverview.Cells.Item(4, 5).Value = "In Translation"
$wsOverview.Cells.Item(4, 6).Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(2, 3).Value = "In Translation"
$wsZhCn.Cells.Item(3, 3).Value = "In Translation"
$wsZhCn.Cells.Item(4, 3).Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(2, 3).Value = "In Translation"
$wsDeDe.Cells.Item(3, 3).Value = "In Translation"
$wsDeDe.Cells.Item(4, 3).Value = "In Translation"

# --- Column width changes ---
# Overview!E:F and the "Status" column (C) on zh-cn / de-de get narrower.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
